$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 1.25

$ws.Range("D3").Value = 1.41
$ws.Range("F3").Value = 1.21

$ws.Range("B4").Value = 1.49
$ws.Range("C4").Value = 1.43
$ws.Range("F4").Value = 1.07
$ws.Range("G4").Value = 1

$ws.Range("C5").Value = 1.37
$ws.Range("D5").Value = 1.33
$ws.Range("F5").Value = 1.03

$ws.Range("C6").Value = 1.47
$ws.Range("D6").Value = 1.54
$ws.Range("G6").Value = 1.03

$ws.Range("D7").Value = 1.71
$ws.Range("F7").Value = 1.47
